$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 297
$ws.Range("B2").Value = '''2025'
$ws.Range("C2").Value = '''3579'
$ws.Range("D2").Value = 'Nguyễn Văn F'
$ws.Range("E2").Value = 'Senior'
$ws.Range("F2").Value = 'Talent Acquisition'
$ws.Range("G2").Value = 'Human Resources Division'
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 4
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 5
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 4
$ws.Range("AA2").Value = 4
$ws.Range("AB2").Value = 4
$ws.Range("AC2").Value = 4
$ws.Range("AD2").Value = 2
$ws.Range("AE2").Value = 3
$ws.Range("AF2").Value = 3
$ws.Range("AG2").Value = 4
$ws.Range("AH2").Value = 'Medium'
$ws.Range("AI2").Value = 'Medium'
$ws.Range("AJ2").Value = '2025-11-05 03:56:54'

# Row 3
$ws.Range("A3").Value = 296
$ws.Range("B3").Value = '''2025'
$ws.Range("C3").Value = '''8911'
$ws.Range("D3").Value = 'Nguyễn Văn E'
$ws.Range("E3").Value = 'Senior'
$ws.Range("F3").Value = 'Brokerage Management Department'
$ws.Range("G3").Value = 'Brokerage Division'
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 4
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 4
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = 5
$ws.Range("U3").Value = 5
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = 2
$ws.Range("X3").Value = 3
$ws.Range("Y3").Value = 2
$ws.Range("Z3").Value = 4
$ws.Range("AA3").Value = 3
$ws.Range("AB3").Value = 3
$ws.Range("AC3").Value = 3
$ws.Range("AD3").Value = 4
$ws.Range("AE3").Value = 5
$ws.Range("AF3").Value = 5
$ws.Range("AG3").Value = 3
$ws.Range("AH3").Value = 'Low'
$ws.Range("AI3").Value = 'Low'
$ws.Range("AJ3").Value = '2025-11-05 03:56:54'

# Row 4
$ws.Range("A4").Value = 295
$ws.Range("B4").Value = '''2025'
$ws.Range("C4").Value = '''4567'
$ws.Range("D4").Value = 'Nguyễn Văn D'
$ws.Range("E4").Value = 'Officer'
$ws.Range("F4").Value = 'Brokerage Management Department'
$ws.Range("G4").Value = 'Brokerage Division'
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 4
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = 2
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 5
$ws.Range("X4").Value = 1
$ws.Range("Y4").Value = 1
$ws.Range("Z4").Value = 4
$ws.Range("AA4").Value = 4
$ws.Range("AB4").Value = 4
$ws.Range("AC4").Value = 4
$ws.Range("AD4").Value = 5
$ws.Range("AE4").Value = 3
$ws.Range("AF4").Value = 3
$ws.Range("AG4").Value = 4
$ws.Range("AH4").Value = 'Low'
$ws.Range("AI4").Value = 'Medium'
$ws.Range("AJ4").Value = '2025-11-05 03:56:54'

# Row 5
$ws.Range("A5").Value = 294
$ws.Range("B5").Value = '''2025'
$ws.Range("C5").Value = '''5678'
$ws.Range("D5").Value = 'Nguyễn Văn C'
$ws.Range("E5").Value = 'Senior'
$ws.Range("F5").Value = 'Talent Acquisition'
$ws.Range("G5").Value = 'Human Resources Division'
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = 4
$ws.Range("O5").Value = 5
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 4
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = 4
$ws.Range("V5").Value = 4
$ws.Range("W5").Value = 3
$ws.Range("X5").Value = 5
$ws.Range("Y5").Value = 3
$ws.Range("Z5").Value = 4
$ws.Range("AA5").Value = 4
$ws.Range("AB5").Value = 4
$ws.Range("AC5").Value = 4
$ws.Range("AD5").Value = 2
$ws.Range("AE5").Value = 3
$ws.Range("AF5").Value = 3
$ws.Range("AG5").Value = 4
$ws.Range("AH5").Value = 'Medium'
$ws.Range("AI5").Value = 'Medium'
$ws.Range("AJ5").Value = '2025-11-05 03:55:04'

# Row 6
$ws.Range("A6").Value = 293
$ws.Range("B6").Value = '''2025'
$ws.Range("C6").Value = '''1345'
$ws.Range("D6").Value = 'Nguyễn Văn B'
$ws.Range("E6").Value = 'Senior'
$ws.Range("F6").Value = 'Brokerage Management Department'
$ws.Range("G6").Value = 'IT Division'
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 2
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 3
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 4
$ws.Range("S6").Value = 5
$ws.Range("T6").Value = 5
$ws.Range("U6").Value = 5
$ws.Range("V6").Value = 1
$ws.Range("W6").Value = 2
$ws.Range("X6").Value = 3
$ws.Range("Y6").Value = 2
$ws.Range("Z6").Value = 4
$ws.Range("AA6").Value = 3
$ws.Range("AB6").Value = 3
$ws.Range("AC6").Value = 3
$ws.Range("AD6").Value = 4
$ws.Range("AE6").Value = 5
$ws.Range("AF6").Value = 5
$ws.Range("AG6").Value = 3
$ws.Range("AH6").Value = 'Low'
$ws.Range("AI6").Value = 'Low'
$ws.Range("AJ6").Value = '2025-11-05 03:55:04'

# Row 7
$ws.Range("A7").Value = 292
$ws.Range("B7").Value = '''2025'
$ws.Range("C7").Value = '''1234'
$ws.Range("D7").Value = 'Nguyễn Văn A'
$ws.Range("E7").Value = 'Officer'
$ws.Range("F7").Value = 'Brokerage Management Department'
$ws.Range("G7").Value = 'Accounting Division'
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 4
$ws.Range("Q7").Value = 4
$ws.Range("R7").Value = 4
$ws.Range("S7").Value = 2
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1
$ws.Range("W7").Value = 5
$ws.Range("X7").Value = 1
$ws.Range("Y7").Value = 1
$ws.Range("Z7").Value = 4
$ws.Range("AA7").Value = 4
$ws.Range("AB7").Value = 4
$ws.Range("AC7").Value = 4
$ws.Range("AD7").Value = 5
$ws.Range("AE7").Value = 3
$ws.Range("AF7").Value = 3
$ws.Range("AG7").Value = 4
$ws.Range("AH7").Value = 'Low'
$ws.Range("AI7").Value = 'Medium'
$ws.Range("AJ7").Value = '2025-11-05 03:55:04'
